# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "K" (strikeout) column values (column G) for the 2021 richards_trevor
# save_data sheet, recomputed from the new K-based calculation.
$updates = @{
    2  = 1
    3  = 0
    5  = 0
    6  = 0
    7  = 3
    9  = 0
    10 = 2
    11 = 2
    12 = 2
    13 = 1
    15 = 1
    16 = 0
    17 = 1
    18 = 0
    19 = 2
    20 = 1
    21 = 1
    22 = 1
    23 = 2
    24 = 0
    25 = 0
    26 = 1
    27 = 2
    28 = 1
    29 = 3
    30 = 2
    31 = 1
    32 = 4
    33 = 1
    34 = 1
    35 = 3
    36 = 1
    37 = 1
    38 = 3
    39 = 1
    40 = 2
    41 = 2
    42 = 4
    43 = 3
    44 = 0
    45 = 0
    46 = 1
    47 = 3
    48 = 1
    49 = 2
    50 = 5
    51 = 0
    52 = 4
    53 = 3
    54 = 3
    55 = 3
    56 = 4
    57 = 1
}

foreach ($row in $updates.Keys) {
    $ws.Range("G$row").Value = $updates[$row]
}
